# Update average_county_temperature (col I), worst_ashp_cop (col N),
# and best_ashp_cop (col O) for the facility rows whose NOAA-derived
# county temperature changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (I, N, O)
$updates = @{
    2  = @(-3.222222222222223, 1.475542118432027, 1.575)
    3  = @(0.4166666666666667, 1.516977491961415, 1.622926829268293)
    4  = @(12.51681286549706,  1.673218141204726, 1.805631235675374)
    7  = @(12.51681286549706,  1.673218141204726, 1.805631235675374)
    9  = @(12.51681286549706,  1.673218141204726, 1.805631235675374)
    10 = @(0.4166666666666667, 1.516977491961415, 1.622926829268293)
    13 = @(15.74228395061728,  1.720452734369724, 1.861492917301914)
    15 = @(15.74228395061728,  1.720452734369724, 1.861492917301914)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value  = $vals[0]  # column I
    $ws.Cells.Item($row, 14).Value = $vals[1]  # column N
    $ws.Cells.Item($row, 15).Value = $vals[2]  # column O
}
